$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.692.10"
$ws.Range("E2").Value = "  -5.05%  "
$ws.Range("D3").Value = "3.069.09"
$ws.Range("E3").Value = "  -5.26%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'539.74"
$ws.Range("E5").Value = "  -7.13%  "
$ws.Range("D6").Value = "'133.71"
$ws.Range("E6").Value = "  -11.86%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "3.062.90"
$ws.Range("E8").Value = "  -5.17%  "
$ws.Range("D9").Value = "'0.489"
$ws.Range("E9").Value = "  -4.77%  "
$ws.Range("D10").Value = "'0.154"
$ws.Range("E10").Value = "  -5.26%  "
$ws.Range("D11").Value = "'6.21"
$ws.Range("E11").Value = "  -12.69%  "
$ws.Range("D12").Value = "'0.461"
$ws.Range("E12").Value = "  -5.35%  "
$ws.Range("D13").Value = "'34.68"
$ws.Range("E13").Value = "  -7.87%  "
$ws.Range("D14").Value = "'0.0000221"
$ws.Range("E14").Value = "  -5.62%  "
$ws.Range("D15").Value = "3.522.60"
$ws.Range("E15").Value = "  -6.08%  "
$ws.Range("D16").Value = "62.560.56"
$ws.Range("E16").Value = "  -5.36%  "
$ws.Range("E17").Value = "  -3.03%  "
$ws.Range("D18").Value = "3.059.63"
$ws.Range("E18").Value = "  -5.56%  "
$ws.Range("D19").Value = "'6.63"
$ws.Range("E19").Value = "  -6.61%  "
$ws.Range("D20").Value = "'480.88"
$ws.Range("E20").Value = "  -11.79%  "
$ws.Range("D21").Value = "'13.37"
$ws.Range("E21").Value = "  -7.79%  "
$ws.Range("D22").Value = "'0.705"
$ws.Range("E22").Value = "  -5.16%  "
$ws.Range("D23").Value = "'7.23"
$ws.Range("E23").Value = "  -8.25%  "
$ws.Range("D24").Value = "'78.56"
$ws.Range("E24").Value = "  -3.04%  "
$ws.Range("D25").Value = "'12.05"
$ws.Range("E25").Value = "  -10.44%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("E27").Value = "  -9.11%  "
$ws.Range("D28").Value = "'8.18"
$ws.Range("E28").Value = "  -11.80%  "
$ws.Range("D29").Value = "'0.998"
$ws.Range("E29").Value = "  -0.23%  "
$ws.Range("D30").Value = "'1.92"
$ws.Range("E30").Value = "  -14.56%  "
$ws.Range("D31").Value = "'26.05"
$ws.Range("E31").Value = "  -5.79%  "
$ws.Range("D32").Value = "'1.10"
$ws.Range("E32").Value = "  -6.57%  "
$ws.Range("B33").Value = "OKB"
$ws.Range("C33").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D33").Value = "'58.87"
$ws.Range("E33").Value = "  +6.82%  "
$ws.Range("B34").Value = "Stacks"
$ws.Range("C34").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D34").Value = "'2.41"
$ws.Range("E34").Value = "  -12.23%  "
$ws.Range("D35").Value = "'6.00"
$ws.Range("E35").Value = "  -5.29%  "
$ws.Range("D36").Value = "'481.67"
$ws.Range("E36").Value = "  -15.16%  "
$ws.Range("D37").Value = "'5.15"
$ws.Range("E37").Value = "  -8.83%  "
$ws.Range("D38").Value = "3.128.89"
$ws.Range("E38").Value = "  -2.00%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "'0.0796"
$ws.Range("E39").Value = "  -7.50%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "'0.0390"
$ws.Range("E40").Value = "  -13.66%  "
$ws.Range("E41").Value = "  -10.24%  "
$ws.Range("D42").Value = "'8.06"
$ws.Range("E42").Value = "  -6.05%  "
$ws.Range("E43").Value = "  -13.59%  "
$ws.Range("E44").Value = "  -10.36%  "
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("E46").Value = "  -11.24%  "
$ws.Range("D47").Value = "'24.65"
$ws.Range("E47").Value = "  -6.79%  "
$ws.Range("D48").Value = "'118.46"
$ws.Range("E48").Value = "  -6.06%  "
$ws.Range("E49").Value = "  -4.57%  "
$ws.Range("D50").Value = "0.0₃0506"
$ws.Range("E50").Value = "  -9.09%  "
$ws.Range("E51").Value = "  -8.69%  "
